# "Remove 'keywords' from Performance model"
#
# The Performance-model row (row 9) had a duplicated "-" placeholder cell
# in column O (it duplicated the value already present in H9/I9/K9/N9).
# That extra "keywords" column is removed here, which shifts every column
# from P onward one position to the left (P->O, Q->P, R->Q, S->R, T->S,
# U->T) and shrinks the sheet's used range from A1:U9 to A1:T9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Columns("O").Delete()

# Reflect the current selection: the author had column N selected after
# making the edit.
[void]$ws.Columns("N").Select()
